$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- The table is shrinking from 18 data rows (16:33) to 12 (16:27), so the
# new last row (27) must take on the thicker "closing" border that the old
# last row (33) had. Copy that formatting over first, while row 33 still
# exists with its original formatting intact.
$ws.Range("B33:J33").Copy() | Out-Null
$ws.Range("B27:J27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Remove the 6 surplus data rows (old rows 28:33). This shifts the
# trailing blank gap + footer block (old rows 34:39) up by 6 rows, so the
# footer lands on rows 32:33, matching the target layout.
$ws.Rows("28:33").Delete() | Out-Null

# --- Rewrite the account-statement table (rows 16-27) with the new dataset.
$tableData = @(
    @(16, "CC", "20090288",        "CARLOS ANDRES BLANCO TUIRAN",      "1709", 29509, 781242),
    @(17, "CC", "73008944",        "ELKIN RAFAEL SIERRA CARO",         "1802", 28290, 848714),
    @(18, "CC", "79417905",        "JUAN ALBERTO HOYOS CUARTAS",       "1905", 33125, 908526),
    @(19, "CC", "93355255",        "LEONEL TORRES",                    "1912", 33125, 828116),
    @(20, "PE", "963298814051990", "STEPHANYE PATRICIA PERDOMO SAER",  "1912", 37276, 931889),
    @(21, "CC", "1047409424",      "JULIETT PAOLA ANGULO BEJARANO",    "2002", 35112, 877803),
    @(22, "CC", "20255153",        "ROSANGELA CAROLINA ROMERO BURGOS", "2002", 35112, 877803),
    @(23, "CC", "1047409424",      "JULIETT PAOLA ANGULO BEJARANO",    "2003", 35112, 877803),
    @(24, "CC", "20255153",        "ROSANGELA CAROLINA ROMERO BURGOS", "2003", 35112, 877803),
    @(25, "CC", "20255153",        "ROSANGELA CAROLINA ROMERO BURGOS", "2004", 35112, 877803),
    @(26, "CC", "1102839947",      "KARINA MARIA BELLO GOMEZCACERES",  "2102", 13325, 908526),
    @(27, "CC", "1047482816",      "DAMIAN PEREZ ACEVEDO",             "2108", 36341, 1200000)
)

foreach ($row in $tableData) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}

# --- Update the summary figures above the table.
$ws.Range("E11").Value = 386551
$ws.Range("C13").Value = 9
$ws.Range("F13").Value = 9
